$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 looks numeric ("001") so force it to be stored as text, then reset the
# cell style back to Normal so no stray style index is left on the cell.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").Style = "Normal"

$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 55401748.66
$ws.Range("P2").Value = 233048588.85
$ws.Range("Q2").Value = 147987605.47
$ws.Range("R2").Value = 1.0454050709
$ws.Range("S2").Value = 85671024.20999999
$ws.Range("T2").Value = 85671024.20999999
$ws.Range("U2").Value = 17.2540621425
$ws.Range("V2").Value = 30454387.73
$ws.Range("W2").Value = 18514204.16
$ws.Range("X2").Value = -1500690.12
$ws.Range("Y2").Value = 83520681.84999999
$ws.Range("Z2").Value = 83941585.25
$ws.Range("AA2").Value = 12184693.8

$ws.Range("AG2").Value = 2637717.96

$ws.Range("AP2").Value = 13.638008843
$ws.Range("AQ2").Value = 44.280807131286
$ws.Range("AR2").Value = 214.984506523406
$ws.Range("AS2").Value = 47636506.43
$ws.Range("AT2").Value = 4.336453488329
